$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits that do not depend on the row insert -------------------------------
$ws.Range("A6").Value = "Vous trouverez ci-dessous mon temps de travail pour chacune des classes."
$ws.Range("B8").Value = "Heure de travail"

# --- Insert a new blank row for the "JUNIT" entry (pushes old rows 12-16 to 13-17) --
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).RowHeight = $ws.Rows.Item(13).RowHeight()

# --- Fix up cell formatting (borders/alignment) so every cell ends up with the -----
# --- same visual style as in the final workbook, using format-only paste. ----------
$xlPasteFormats = -4122

function CopyFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

CopyFormat "C11" "C14"
CopyFormat "B15" "B16"
CopyFormat "C15" "C16"
CopyFormat "A13" "A10"
CopyFormat "C10" "C11"
CopyFormat "A11" "A12"
CopyFormat "B10" "B12"
CopyFormat "C10" "C12"
CopyFormat "A11" "A13"
CopyFormat "C10" "C13"
CopyFormat "A10" "A15"
CopyFormat "B10" "B15"
CopyFormat "C10" "C15"

$excel.CutCopyMode = $false

# --- Final content --------------------------------------------------------------
$ws.Range("A9").Value = "Calculatrice"
$ws.Range("B9").Value = 25
$ws.Range("C9").Value = "- Interface graphique`n- Algorythme"

$ws.Range("A10").Value = "PanelEcranCenter`nPanelEcranNorth"
$ws.Range("B10").Value = 25
$ws.Range("C10").Value = "- Interface graphique (applications, disposition,…)`n- SSID, signal`n- Heure"

$ws.Range("A11").Value = "Gallery"
$ws.Range("B11").Value = 25
$ws.Range("C11").Value = "- Interface graphique`n- Ajout d'image`n- Suppression d'images"

$ws.Range("A12").Value = "JUNIT"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""

$ws.Range("A13").Value = "Settings"
$ws.Range("B13").Value = 8
$ws.Range("C13").Value = "- Interface graphique`n- Liaison avec la Gallery"

$ws.Range("A14").Value = "Frame"
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = "- Interface graphique"

$ws.Range("A15").Value = "IconButton`nIconPanel"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = "- Exceptions`n- Image draw`n- Méthodes"

$ws.Range("A16").Value = "Javadoc + diverses modifications"
$ws.Range("B16").Value = 4
$ws.Range("C16").Value = "- Ajout des commentaires`n- Optimisation du code"

$ws.Range("A17").Value = "TOTAL"
$ws.Range("B17").Formula = "=SUM(B9:B16)"
$ws.Range("C17").Value = ""

# --- Selection / view matches the saved workbook state ---------------------------
$ws.Range("C12").Select()
